$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the "Missed Annotations" column (old column B) ---
# Shift columns C..H left into B..G using Cut/Insert (preserves exact stored
# column widths instead of the lossy ColumnWidth-character rounding that a
# plain width assignment would introduce).
$ws.Columns(3).Cut()
$ws.Columns(2).Insert()
$ws.Columns(4).Cut()
$ws.Columns(3).Insert()
$ws.Columns(5).Cut()
$ws.Columns(4).Insert()
$ws.Columns(6).Cut()
$ws.Columns(5).Insert()
$ws.Columns(7).Cut()
$ws.Columns(6).Insert()
$ws.Columns(8).Cut()
$ws.Columns(7).Insert()

# Column H now holds the old "Missed Annotations" column (leftover). Clear it;
# it will be rebuilt below as the new "Notes" column.
$ws.Columns(8).ClearContents()

# --- Step 2: column widths for the new layout ---
$ws.Columns(7).ColumnWidth = 30.166666666666668
$ws.Columns(8).ColumnWidth = 61.307291666666664

# --- Step 3: fill in the data column-by-column (this also matches the order
# new strings get interned into the shared-string table) ---
$ws.Range("A4").Value = "SharpChecker"
$ws.Range("A5").Value = "EventHub"

$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 1

$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 0

$ws.Range("D4").Value = 4368
$ws.Range("D5").Value = 268144

$ws.Range("E4").Value = 7
$ws.Range("E5").Value = 9

$ws.Range("F4").Value = "Nullness"
$ws.Range("F5").Value = "Nullness"

$ws.Range("G3").Value = "Explicit Checks/Assertions Added"
$ws.Range("G4").Value = 11
$ws.Range("G5").Value = 0

$ws.Range("H3").Value = "Notes"
$ws.Range("H4").Value = "Null reference exceptions were occurring when analyzing EventHub, and the class hierarchy issue was uncovered as a result of the method override checking (a little stretch)."
$ws.Range("H5").Value = "There is an invocation which should be presenting a diagnostic, but is not.  It may be that that project is not being analzyed properly."

# --- Step 6: formatting ---
$ws.Range("H3:H5").WrapText = $true
$ws.Rows(4).RowHeight = 45
$ws.Rows(5).RowHeight = 30

# --- Step 7: page setup / selection ---
$ws.PageSetup.Orientation = 1
$ws.Range("H5").Select()
